$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in missing unit values with "-" placeholder (added first so this string gets the earliest new index)
$ws.Cells.Item(57,3).Value = "-"
$ws.Cells.Item(58,2).Value = "-"
$ws.Cells.Item(58,3).Value = "-"

# Add new "Category" column (E) header and per-row classification
$ws.Cells.Item(1,5).Value = "Category"
$ws.Cells.Item(2,5).Value = "Meat"
$ws.Cells.Item(3,5).Value = "Meat"
$ws.Cells.Item(4,5).Value = "Meat"
$ws.Cells.Item(5,5).Value = "Meat"
$ws.Cells.Item(6,5).Value = "Meat"
$ws.Cells.Item(7,5).Value = "Meat"
$ws.Cells.Item(8,5).Value = "Meat"
$ws.Cells.Item(9,5).Value = "Meat"
$ws.Cells.Item(10,5).Value = "Meat"
$ws.Cells.Item(11,5).Value = "Meat"
$ws.Cells.Item(12,5).Value = "Meat"
$ws.Cells.Item(13,5).Value = "Dairy"
$ws.Cells.Item(14,5).Value = "Dairy"
$ws.Cells.Item(15,5).Value = "Dairy"
$ws.Cells.Item(16,5).Value = "Dairy"
$ws.Cells.Item(17,5).Value = "Dairy"
$ws.Cells.Item(18,5).Value = "Meat"
$ws.Cells.Item(19,5).Value = "Grain"
$ws.Cells.Item(20,5).Value = "Grain"
$ws.Cells.Item(21,5).Value = "Grain"
$ws.Cells.Item(22,5).Value = "Grain"
$ws.Cells.Item(23,5).Value = "Grain"
$ws.Cells.Item(24,5).Value = "Fruit"
$ws.Cells.Item(25,5).Value = "Fruit"
$ws.Cells.Item(26,5).Value = "Fruit"
$ws.Cells.Item(27,5).Value = "Fruit"
$ws.Cells.Item(28,5).Value = "Fruit"
$ws.Cells.Item(29,5).Value = "Fruit"
$ws.Cells.Item(30,5).Value = "Vegetable"
$ws.Cells.Item(31,5).Value = "Vegetable"
$ws.Cells.Item(32,5).Value = "Other"
$ws.Cells.Item(33,5).Value = "Vegetable"
$ws.Cells.Item(34,5).Value = "Vegetable"
$ws.Cells.Item(35,5).Value = "Vegetable"
$ws.Cells.Item(36,5).Value = "Vegetable"
$ws.Cells.Item(37,5).Value = "Vegetable"
$ws.Cells.Item(38,5).Value = "Vegetable"
$ws.Cells.Item(39,5).Value = "Other"
$ws.Cells.Item(40,5).Value = "Other"
$ws.Cells.Item(41,5).Value = "Other"
$ws.Cells.Item(42,5).Value = "Other"
$ws.Cells.Item(43,5).Value = "Other"
$ws.Cells.Item(44,5).Value = "Other"
$ws.Cells.Item(45,5).Value = "Other"
$ws.Cells.Item(46,5).Value = "Other"
$ws.Cells.Item(47,5).Value = "Other"
$ws.Cells.Item(48,5).Value = "Other"
$ws.Cells.Item(49,5).Value = "Other"
$ws.Cells.Item(50,5).Value = "Other"
$ws.Cells.Item(51,5).Value = "Other"
$ws.Cells.Item(52,5).Value = "Other"
$ws.Cells.Item(53,5).Value = "Other"
$ws.Cells.Item(54,5).Value = "Other"
$ws.Cells.Item(55,5).Value = "Other"
$ws.Cells.Item(56,5).Value = "Other"
$ws.Cells.Item(57,5).Value = "Other"
$ws.Cells.Item(58,5).Value = "Other"

# Scroll the view and move the active selection, as left by the author
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("E59").Select()
